$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 490; this shifts the existing rows
# 490-553 down to become rows 492-555.
$ws.Rows.Item(490).Insert()
$ws.Rows.Item(490).Insert()

# Populate the newly inserted row 490 with its new record.
$ws.Range("A490").Value = 6
$ws.Range("B490").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C490").Value = "Metropolitana"
$ws.Range("D490").Value = 44491
$ws.Range("E490").Value = 13
$ws.Range("F490").Value = 100112040
$ws.Range("G490").Value = "Cilantro"
$ws.Range("H490").Value = "Sin especificar"
$ws.Range("I490").Value = "Primera"
$ws.Range("J490").Value = 680
$ws.Range("K490").Value = 4000
$ws.Range("L490").Value = 4500
$ws.Range("M490").Value = 4235
$ws.Range("N490").Value = "`$/caja 36 atados"
$ws.Range("O490").Value = "Región Metropolitana"
$ws.Range("P490").Value = 118
$ws.Range("Q490").Value = 36
$ws.Range("R490").Value = "Hortaliza"

# Populate the newly inserted row 491 with its new record.
$ws.Range("A491").Value = 6
$ws.Range("B491").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C491").Value = "Metropolitana"
$ws.Range("D491").Value = 44491
$ws.Range("E491").Value = 13
$ws.Range("F491").Value = 100112040
$ws.Range("G491").Value = "Cilantro"
$ws.Range("H491").Value = "Sin especificar"
$ws.Range("I491").Value = "Primera"
$ws.Range("J491").Value = 470
$ws.Range("K491").Value = 7000
$ws.Range("L491").Value = 7500
$ws.Range("M491").Value = 7202
$ws.Range("N491").Value = "`$/docena de atados"
$ws.Range("O491").Value = "Región Metropolitana"
$ws.Range("P491").Value = 2401
$ws.Range("Q491").Value = 3
$ws.Range("R491").Value = "Hortaliza"
